$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers stay the same (row 1). We only need to rebuild the data rows
# (rows 2-7) to reflect the new combinations of Sending cluster / Target cluster.

$data = @(
    @("FAPs","Fgf17","Fgfr4","ECs", 1, 0.3333333333333333, 0.05925133333333333, 0.177754, 0.129311749608618, 0.129311749608618, 1, 0.3333333333333333, 0.5555316666666666, 1.666595, 0.01938483203642842, 0.01938483203642843, 0.03291599195888888, 0.29624392763, 0.002506686546499748, 0.002506686546499749),
    @("FAPs","Fgf17","Fgfr4","FAPs", 1, 0.3333333333333333, 0.05925133333333333, 0.177754, 0.129311749608618, 0.129311749608618, 2, 0.6666666666666666, 0.1817723333333333, 0.5453170000000001, 0.006342799811357313, 0.006342799811357313, 0.01077025311311111, 0.096932278018, 0.0008201985410238261, 0.0008201985410238262),
    @("FAPs","Fgf17","Fgfr4","sCs", 1, 0.3333333333333333, 0.05925133333333333, 0.177754, 0.129311749608618, 0.129311749608618, 3, 1, 27.92075533333333, 83.762266, 0.9742723681522142, 0.9742723681522143, 1.654341981173778, 14.889077830564, 0.1259848645210944, 0.1259848645210944),
    @("sCs","Fgf17","Fgfr4","ECs", 3, 1, 0.398954, 1.196862, 0.8706882503913821, 0.8706882503913821, 1, 0.3333333333333333, 0.5555316666666666, 1.666595, 0.01938483203642842, 0.01938483203642843, 0.2216315805433333, 1.99468422489, 0.01687814548992868, 0.01687814548992868),
    @("sCs","Fgf17","Fgfr4","FAPs", 3, 1, 0.398954, 1.196862, 0.8706882503913821, 0.8706882503913821, 2, 0.6666666666666666, 0.1817723333333333, 0.5453170000000001, 0.006342799811357313, 0.006342799811357313, 0.07251879947266668, 0.6526691952540001, 0.005522601270333487, 0.005522601270333488),
    @("sCs","Fgf17","Fgfr4","sCs", 3, 1, 0.398954, 1.196862, 0.8706882503913821, 0.8706882503913821, 3, 1, 27.92075533333333, 83.762266, 0.9742723681522142, 0.9742723681522143, 11.13909702325467, 100.251873209292, 0.8482875036311198, 0.8482875036311199)
)

$row = 2
foreach ($rec in $data) {
    $col = 1
    foreach ($val in $rec) {
        $ws.Cells.Item($row, $col).Value = $val
        $col++
    }
    $row++
}
